$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary"
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("F2").Value = 0
$wsSummary.Range("A3").Value = 700
$wsSummary.Range("E3").Value = 500

# ---------------------------------------------------------------------------
# Sheet "Repayment schedule"
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Row 4 - shifts up from what used to be row 5's timing, and amounts recalc
$wsRepay.Range("B4").Value = 31
$wsRepay.Range("C4").Value = 42095
$wsRepay.Range("H4").Value = 90.91

# K4 / P4 become round numbers -> take on the "#,##0" look that K5/P5 used to have
$wsRepay.Range("K5").Copy()
$wsRepay.Range("K4").PasteSpecial(-4122)
$wsRepay.Range("K4").Value = 1000

$wsRepay.Range("P5").Copy()
$wsRepay.Range("P4").PasteSpecial(-4122)
$wsRepay.Range("P4").Value = 1000

# Row 5
$wsRepay.Range("B5").Value = 30
$wsRepay.Range("C5").Value = 42125
$wsRepay.Range("H5").Value = 81.82

# K5 / P5 become plain decimals -> take on the "General" look from L5/L4 etc.
$wsRepay.Range("L5").Copy()
$wsRepay.Range("K5").PasteSpecial(-4122)
$wsRepay.Range("K5").Value = 990.91

$wsRepay.Range("L4").Copy()
$wsRepay.Range("P5").PasteSpecial(-4122)
$wsRepay.Range("P5").Value = 990.91

# Row 6
$wsRepay.Range("B6").Value = 31
$wsRepay.Range("C6").Value = 42156

# Row 7
$wsRepay.Range("B7").Value = 30
$wsRepay.Range("C7").Value = 42186

# Row 8
$wsRepay.Range("B8").Value = 31
$wsRepay.Range("C8").Value = 42217

# Row 9 (only the date moves)
$wsRepay.Range("C9").Value = 42248

# Row 10
$wsRepay.Range("B10").Value = 30
$wsRepay.Range("C10").Value = 42278

# Row 11
$wsRepay.Range("B11").Value = 31
$wsRepay.Range("C11").Value = 42309

# Row 12
$wsRepay.Range("B12").Value = 30
$wsRepay.Range("C12").Value = 42339

# Row 13
$wsRepay.Range("B13").Value = 31
$wsRepay.Range("C13").Value = 42370

# The whole "heading" column (O) of the detail rows, plus the stray P2 cell,
# get wiped out completely (contents + formatting) so the exporter drops them.
$wsRepay.Range("P2").Clear()
$wsRepay.Range("O3:O13").Clear()

# ---------------------------------------------------------------------------
# Sheet "Transactions"
# ---------------------------------------------------------------------------
$wsTrans = $wb.Worksheets.Item("Transactions")

$wsTrans.Range("A2").Value = 6372
$wsTrans.Range("A3").Value = 6370

# ---------------------------------------------------------------------------
# Selections / active tab: Summary -> Repayment schedule -> Transactions,
# finishing on Transactions so it ends up the active sheet/tab.
# ---------------------------------------------------------------------------
$wsSummary.Activate()
$wsSummary.Range("E5").Select()

$wsRepay.Activate()
$wsRepay.Range("O11").Select()

$wsTrans.Activate()
$wsTrans.Range("D3").Select()
